# Update "想去人数" (F column) counts on the 展览, 演出 and 全部类型 sheets
# to match the regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 342
$ws1.Range("F3").Value = 113
$ws1.Range("F5").Value = 5043
$ws1.Range("F8").Value = 295
$ws1.Range("F9").Value = 763
$ws1.Range("F10").Value = 246

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 31

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 342
$ws4.Range("F3").Value = 113
$ws4.Range("F5").Value = 5043
$ws4.Range("F8").Value = 295
$ws4.Range("F9").Value = 763
$ws4.Range("F10").Value = 31
$ws4.Range("F11").Value = 246
